# Update cryptos list figures (price + 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.172.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "'2.509.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("D5").Value = "'109.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").Value = "'320.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.547"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.45%  "
$ws.Range("D10").Value = "'39.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.22%  "
$ws.Range("D11").Value = "'20.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.69%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "'2.902.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "'2.513.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "'48.014.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("D19").Value = "'13.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("D23").Value = "'72.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.61%  "
$ws.Range("D24").Value = "'274.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.82%  "
$ws.Range("D25").Value = "'2.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("E28").Value = "  +5.92%  "
$ws.Range("D29").Value = "'10.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").Value = "'0.141"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.24%  "
$ws.Range("D31").Value = "'35.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("D32").Value = "'49.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("D33").Value = "'19.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.54%  "
$ws.Range("D34").Value = "'5.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'0.0782"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "'1.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("D39").Value = "'2.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("D41").Value = "'121.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.23%  "
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").Value = "'21.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.39%  "
$ws.Range("D44").Value = "'0.0308"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.19%  "
$ws.Range("D45").Value = "'2.030.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("E47").Value = "  +5.18%  "
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("D50").Value = "'5.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.82%  "
$ws.Range("D51").Value = "'79.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.34%  "
